# Auto-generated edits applying the profit-sheet numeric updates described in the diff.
# (commit: "chore: update Sheets via scheduled runner")
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 13 (ALC)
$ws.Range("H13").Value = 30000
$ws.Range("J13").Value = 30000
$ws.Range("L13").Value = 30000
$ws.Range("N13").Value = -30338

# Row 94 (ALC)
$ws.Range("H94").Value = 1300
$ws.Range("I94").Value = 1300
$ws.Range("K94").Value = 1300
$ws.Range("M94").Value = -849

# Row 137 (ALC)
$ws.Range("H137").Value = 825.2857
$ws.Range("I137").Value = 799.2222
$ws.Range("J137").Value = 913.25
$ws.Range("K137").Value = 2397.6666
$ws.Range("L137").Value = 2739.75
$ws.Range("M137").Value = 152.3334
$ws.Range("N137").Value = -7839.75

$ws = $wb.Worksheets.Item("ARM")
# Row 64 (ARM)
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496

# Row 67 (ARM)
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716

# Row 88 (ARM)
$ws.Range("H88").Value = 1724.5454
$ws.Range("I88").Value = 1452.2222
$ws.Range("J88").Value = 2950
$ws.Range("K88").Value = 1452.2222
$ws.Range("L88").Value = 2950
$ws.Range("M88").Value = -1046.2222
$ws.Range("N88").Value = -3762

# Row 91 (ARM)
$ws.Range("H91").Value = 1724.5454
$ws.Range("I91").Value = 1452.2222
$ws.Range("J91").Value = 2950
$ws.Range("K91").Value = 1452.2222
$ws.Range("L91").Value = 2950
$ws.Range("M91").Value = -48.22219999999993
$ws.Range("N91").Value = -5758

# Row 105 (ARM)
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 86 (BSM)
$ws.Range("H86").Value = 1631.2
$ws.Range("I86").Value = 1321.5385
$ws.Range("J86").Value = 1966.6666
$ws.Range("K86").Value = 1321.5385
$ws.Range("L86").Value = 1966.6666
$ws.Range("M86").Value = -198.5385000000001
$ws.Range("N86").Value = -4212.6666

# Row 89 (BSM)
$ws.Range("H89").Value = 1631.2
$ws.Range("I89").Value = 1321.5385
$ws.Range("J89").Value = 1966.6666
$ws.Range("K89").Value = 6607.692500000001
$ws.Range("L89").Value = 9833.333000000001
$ws.Range("M89").Value = -991.692500000001
$ws.Range("N89").Value = -21065.333

# Row 99 (BSM)
$ws.Range("H99").Value = 1748.4286
$ws.Range("I99").Value = 1752.5454
$ws.Range("J99").Value = 1733.3334
$ws.Range("K99").Value = 1752.5454
$ws.Range("L99").Value = 1733.3334
$ws.Range("M99").Value = -254.5454
$ws.Range("N99").Value = -4729.3334

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 2263.7715
$ws.Range("I31").Value = 1909.0769
$ws.Range("J31").Value = 3288.4443
$ws.Range("K31").Value = 1909.0769
$ws.Range("L31").Value = 3288.4443
$ws.Range("M31").Value = -1614.0769
$ws.Range("N31").Value = -3878.4443

# Row 34 (CRP)
$ws.Range("H34").Value = 2263.7715
$ws.Range("I34").Value = 1909.0769
$ws.Range("J34").Value = 3288.4443
$ws.Range("K34").Value = 1909.0769
$ws.Range("L34").Value = 3288.4443
$ws.Range("M34").Value = -1707.0769
$ws.Range("N34").Value = -3692.4443

# Row 62 (CRP)
$ws.Range("H62").Value = 3237.2856
$ws.Range("I62").Value = 3483.1428
$ws.Range("J62").Value = 2991.4285
$ws.Range("K62").Value = 3483.1428
$ws.Range("L62").Value = 2991.4285
$ws.Range("M62").Value = -2859.1428
$ws.Range("N62").Value = -4239.4285

# Row 65 (CRP)
$ws.Range("H65").Value = 3237.2856
$ws.Range("I65").Value = 3483.1428
$ws.Range("J65").Value = 2991.4285
$ws.Range("K65").Value = 17415.714
$ws.Range("L65").Value = 14957.1425
$ws.Range("M65").Value = -14295.714
$ws.Range("N65").Value = -21197.1425

$ws = $wb.Worksheets.Item("CUL")
# Row 101 (CUL)
$ws.Range("H101").Value = 8585.714
$ws.Range("J101").Value = 8585.714
$ws.Range("L101").Value = 25757.142
$ws.Range("N101").Value = -30625.142

# Row 131 (CUL)
$ws.Range("H131").Value = 471.2
$ws.Range("I131").Value = 298.375
$ws.Range("J131").Value = 1162.5
$ws.Range("K131").Value = 895.125
$ws.Range("L131").Value = 3487.5
$ws.Range("M131").Value = 4144.875
$ws.Range("N131").Value = -13567.5

# Row 132 (CUL)
$ws.Range("H132").Value = 1280.6666
$ws.Range("I132").Value = 1146
$ws.Range("J132").Value = 1550
$ws.Range("K132").Value = 10314
$ws.Range("L132").Value = 13950
$ws.Range("M132").Value = -7784
$ws.Range("N132").Value = -19010

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (GSM)
$ws.Range("H2").Value = 53.25
$ws.Range("I2").Value = 46.57143
$ws.Range("K2").Value = 46.57143
$ws.Range("M2").Value = 66.42857000000001

# Row 70 (GSM)
$ws.Range("H70").Value = 5911.4
$ws.Range("I70").Value = 5292.1763
$ws.Range("J70").Value = 7227.25
$ws.Range("K70").Value = 5292.1763
$ws.Range("L70").Value = 7227.25
$ws.Range("M70").Value = -5022.1763
$ws.Range("N70").Value = -7767.25

# Row 73 (GSM)
$ws.Range("H73").Value = 5911.4
$ws.Range("I73").Value = 5292.1763
$ws.Range("J73").Value = 7227.25
$ws.Range("K73").Value = 5292.1763
$ws.Range("L73").Value = 7227.25
$ws.Range("M73").Value = -4356.1763
$ws.Range("N73").Value = -9099.25

# Row 97 (GSM)
$ws.Range("H97").Value = 456.73334
$ws.Range("I97").Value = 334.6154
$ws.Range("J97").Value = 1250.5
$ws.Range("K97").Value = 334.6154
$ws.Range("L97").Value = 1250.5
$ws.Range("M97").Value = 161.3846
$ws.Range("N97").Value = -2242.5

$ws = $wb.Worksheets.Item("LTW")
# Row 46 (LTW)
$ws.Range("H46").Value = 1077.8667
$ws.Range("I46").Value = 918.7
$ws.Range("J46").Value = 1396.2
$ws.Range("K46").Value = 918.7
$ws.Range("L46").Value = 1396.2
$ws.Range("M46").Value = -730.7
$ws.Range("N46").Value = -1772.2

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (WVR)
$ws.Range("H81").Value = 7156.35
$ws.Range("I81").Value = 11863.1
$ws.Range("J81").Value = 2449.6
$ws.Range("K81").Value = 23726.2
$ws.Range("L81").Value = 4899.2
$ws.Range("M81").Value = -22665.2
$ws.Range("N81").Value = -7021.2

# Row 84 (WVR)
$ws.Range("H84").Value = 7156.35
$ws.Range("I84").Value = 11863.1
$ws.Range("J84").Value = 2449.6
$ws.Range("K84").Value = 118631
$ws.Range("L84").Value = 24496
$ws.Range("M84").Value = -113327
$ws.Range("N84").Value = -35104

# Row 122 (WVR)
$ws.Range("H122").Value = 1992.2307
$ws.Range("I122").Value = 1840.8
$ws.Range("J122").Value = 2086.875
$ws.Range("K122").Value = 5522.4
$ws.Range("L122").Value = 6260.625
$ws.Range("M122").Value = -3072.4
$ws.Range("N122").Value = -11160.625
